$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 55
$ws.Range("H55").Value = 362
$ws.Range("I55").Value = 554
$ws.Range("J55").Value = 314
$ws.Range("K55").Value = 554
$ws.Range("L55").Value = 314
$ws.Range("M55").Value = -340
$ws.Range("N55").Value = -742
# Row 64
$ws.Range("H64").Value = 10666.25
$ws.Range("I64").Value = 5427.857
$ws.Range("J64").Value = 18000
$ws.Range("K64").Value = 5427.857
$ws.Range("L64").Value = 18000
$ws.Range("M64").Value = -5179.857
$ws.Range("N64").Value = -18496
# Row 67
$ws.Range("H67").Value = 10666.25
$ws.Range("I67").Value = 5427.857
$ws.Range("J67").Value = 18000
$ws.Range("K67").Value = 5427.857
$ws.Range("L67").Value = 18000
$ws.Range("M67").Value = -4569.857
$ws.Range("N67").Value = -19716
# Row 111
$ws.Range("H111").Value = 1351.2354
$ws.Range("I111").Value = 1197.4546
$ws.Range("J111").Value = 1633.1666
$ws.Range("K111").Value = 3592.3638
$ws.Range("L111").Value = 4899.4998
$ws.Range("M111").Value = -525.3638000000001
$ws.Range("N111").Value = -11033.4998
# Row 132
$ws.Range("H132").Value = 3633.9688
$ws.Range("I132").Value = 2474.5
$ws.Range("J132").Value = 5566.4165
$ws.Range("K132").Value = 7423.5
$ws.Range("L132").Value = 16699.2495
$ws.Range("M132").Value = -4893.5
# Row 138
$ws.Range("H138").Value = 4132.364
$ws.Range("I138").Value = 3544.5
$ws.Range("J138").Value = 5700
$ws.Range("K138").Value = 10633.5
$ws.Range("L138").Value = 17100
$ws.Range("M138").Value = -5493.5
$ws.Range("N138").Value = -27380

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 26
$ws.Range("H26").Value = 381.75
$ws.Range("I26").Value = 405.54544
$ws.Range("J26").Value = 120
$ws.Range("K26").Value = 405.54544
$ws.Range("L26").Value = 120
$ws.Range("M26").Value = -75.54543999999999
$ws.Range("N26").Value = -780
# Row 45
$ws.Range("H45").Value = 4954.625
$ws.Range("I45").Value = 1644
$ws.Range("J45").Value = 5427.5713
$ws.Range("K45").Value = 1644
$ws.Range("L45").Value = 5427.5713
$ws.Range("M45").Value = -1267
$ws.Range("N45").Value = -6181.5713
# Row 61
$ws.Range("H61").Value = 971.1111
$ws.Range("I61").Value = 971.1111
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 971.1111
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -759.1111
# Row 74
$ws.Range("H74").Value = 475
$ws.Range("I74").Value = 475
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 475
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = 399
# Row 77
$ws.Range("H77").Value = 475
$ws.Range("I77").Value = 475
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 2375
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = 1993
# Row 88
$ws.Range("H88").Value = 4280.7334
$ws.Range("I88").Value = 3451
$ws.Range("J88").Value = 4582.4546
$ws.Range("K88").Value = 3451
$ws.Range("L88").Value = 4582.4546
$ws.Range("M88").Value = -3045
$ws.Range("N88").Value = -5394.4546
# Row 91
$ws.Range("H91").Value = 4280.7334
$ws.Range("I91").Value = 3451
$ws.Range("J91").Value = 4582.4546
$ws.Range("K91").Value = 3451
$ws.Range("L91").Value = 4582.4546
$ws.Range("M91").Value = -2047
$ws.Range("N91").Value = -7390.4546
# Row 97
$ws.Range("H97").Value = 2470
$ws.Range("I97").Value = 410
$ws.Range("J97").Value = 3500
$ws.Range("K97").Value = 410
$ws.Range("L97").Value = 3500
$ws.Range("M97").Value = 86
$ws.Range("N97").Value = -4492
# Row 136
$ws.Range("H136").Value = 971.1111
$ws.Range("I136").Value = 971.1111
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 2913.3333
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -363.3332999999998

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 3962
$ws.Range("I20").Value = 2109.6667
$ws.Range("J20").Value = 7666.6665
$ws.Range("K20").Value = 2109.6667
$ws.Range("L20").Value = 7666.6665
$ws.Range("M20").Value = -1862.6667
# Row 86
$ws.Range("H86").Value = 3823.96
$ws.Range("I86").Value = 3191.611
$ws.Range("J86").Value = 5450
$ws.Range("K86").Value = 3191.611
$ws.Range("L86").Value = 5450
$ws.Range("M86").Value = -2068.611
# Row 89
$ws.Range("H89").Value = 3823.96
$ws.Range("I89").Value = 3191.611
$ws.Range("J89").Value = 5450
$ws.Range("K89").Value = 15958.055
$ws.Range("L89").Value = 27250
$ws.Range("M89").Value = -10342.055
# Row 94
$ws.Range("H94").Value = 3448.1667
$ws.Range("I94").Value = 344.5
$ws.Range("J94").Value = 5000
$ws.Range("K94").Value = 344.5
$ws.Range("L94").Value = 5000
$ws.Range("M94").Value = 106.5
$ws.Range("N94").Value = -5902

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 1742.8096
$ws.Range("I7").Value = 917.5714
$ws.Range("J7").Value = 3393.2856
$ws.Range("K7").Value = 917.5714
$ws.Range("L7").Value = 3393.2856
$ws.Range("M7").Value = -804.5714
$ws.Range("N7").Value = -3619.2856
# Row 14
$ws.Range("H14").Value = 2000
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 2000
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 2000
$ws.Range("N14").Value = -2340
# Row 31
$ws.Range("H31").Value = 1897.9584
$ws.Range("I31").Value = 1897.9584
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 1897.9584
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -1602.9584
# Row 32
$ws.Range("H32").Value = 1286.6666
$ws.Range("I32").Value = 680
$ws.Range("J32").Value = 2500
$ws.Range("K32").Value = 680
$ws.Range("L32").Value = 2500
$ws.Range("M32").Value = -364
$ws.Range("N32").Value = -3132
# Row 34
$ws.Range("H34").Value = 1897.9584
$ws.Range("I34").Value = 1897.9584
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 1897.9584
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -1695.9584
# Row 47
$ws.Range("H47").Value = 10000
$ws.Range("I47").Value = 10000
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 10000
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = -9434
# Row 63
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("N63").ClearContents()
# Row 66
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("N66").ClearContents()
# Row 132
$ws.Range("H132").Value = 3071.9
$ws.Range("I132").Value = 3452.5
$ws.Range("J132").Value = 1549.5
$ws.Range("K132").Value = 10357.5
$ws.Range("L132").Value = 4648.5
$ws.Range("M132").Value = -7827.5
# Row 134
$ws.Range("H134").Value = 1871.8572
$ws.Range("I134").Value = 1364.6666
$ws.Range("J134").Value = 2252.25
$ws.Range("K134").Value = 4093.9998
$ws.Range("L134").Value = 6756.75
$ws.Range("M134").Value = -1558.9998

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 333333860
$ws.Range("I4").Value = 285714850
$ws.Range("J4").Value = 500000400
$ws.Range("K4").Value = 857144550
$ws.Range("L4").Value = 1500001200
$ws.Range("M4").Value = -857144438
# Row 106
$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("M106").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 51
$ws.Range("H51").Value = 2500
$ws.Range("I51").Value = 2500
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 2500
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -1991
# Row 113
$ws.Range("H113").Value = 1650
$ws.Range("I113").Value = 1650
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1650
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 520
# Row 122
$ws.Range("H122").Value = 1039.2222
$ws.Range("I122").Value = 1044.25
$ws.Range("J122").Value = 999
$ws.Range("K122").Value = 3132.75
$ws.Range("L122").Value = 2997
$ws.Range("M122").Value = -682.75
$ws.Range("N122").Value = -7897

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 11
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("M11").ClearContents()
$ws.Range("N11").ClearContents()
# Row 22
$ws.Range("H22").Value = 794
$ws.Range("I22").Value = 195
$ws.Range("J22").Value = 1992
$ws.Range("K22").Value = 195
$ws.Range("L22").Value = 1992
$ws.Range("M22").Value = 100
# Row 27
$ws.Range("H27").Value = 794
$ws.Range("I27").Value = 195
$ws.Range("J27").Value = 1992
$ws.Range("K27").Value = 195
$ws.Range("L27").Value = 1992
$ws.Range("M27").Value = -88
# Row 55
$ws.Range("H55").Value = 394.13333
$ws.Range("I55").Value = 88.25
$ws.Range("J55").Value = 505.36365
$ws.Range("K55").Value = 88.25
$ws.Range("L55").Value = 505.36365
$ws.Range("M55").Value = 84.75
$ws.Range("N55").Value = -851.36365
# Row 68
$ws.Range("H68").Value = 2266.1667
$ws.Range("I68").Value = 1398.75
$ws.Range("J68").Value = 4001
$ws.Range("K68").Value = 1398.75
$ws.Range("L68").Value = 4001
$ws.Range("M68").Value = -649.75
$ws.Range("N68").Value = -5499
# Row 71
$ws.Range("H71").Value = 2266.1667
$ws.Range("I71").Value = 1398.75
$ws.Range("J71").Value = 4001
$ws.Range("K71").Value = 6993.75
$ws.Range("L71").Value = 20005
$ws.Range("M71").Value = -3249.75
$ws.Range("N71").Value = -27493
# Row 93
$ws.Range("H93").Value = 1434.5834
$ws.Range("I93").Value = 1455.909
$ws.Range("J93").Value = 1200
$ws.Range("K93").Value = 1455.909
$ws.Range("L93").Value = 1200
$ws.Range("M93").Value = -207.9090000000001
$ws.Range("N93").Value = -3696
# Row 136
$ws.Range("H136").Value = 3793.125
$ws.Range("I136").Value = 3057.5
$ws.Range("J136").Value = 6000
$ws.Range("K136").Value = 9172.5
$ws.Range("L136").Value = 18000
$ws.Range("M136").Value = -6622.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 4357.143
$ws.Range("I62").Value = 2750
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 2750
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -2126
# Row 65
$ws.Range("H65").Value = 4357.143
$ws.Range("I65").Value = 2750
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 13750
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -10630
# Row 96
$ws.Range("H96").Value = 3062.375
$ws.Range("I96").Value = 2949.75
$ws.Range("J96").Value = 3175
$ws.Range("K96").Value = 2949.75
$ws.Range("L96").Value = 3175
$ws.Range("M96").Value = -1576.75
$ws.Range("N96").Value = -5921
# Row 113
$ws.Range("H113").Value = 347.25
$ws.Range("I113").Value = 355.5
$ws.Range("J113").Value = 322.5
$ws.Range("K113").Value = 1066.5
$ws.Range("L113").Value = 967.5
$ws.Range("M113").Value = 1103.5
# Row 132
$ws.Range("H132").Value = 1217.7142
$ws.Range("I132").Value = 1228.2941
$ws.Range("J132").Value = 1172.75
$ws.Range("K132").Value = 3684.8823
$ws.Range("L132").Value = 3518.25
$ws.Range("M132").Value = -1154.8823
$ws.Range("N132").Value = -8578.25
